# Add a new "2022" column (S) to the transport statistics table, mirroring
# the formatting of the existing "2021" column (R), and update the active
# selection to T3 (the cell to the right of the new column's header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New data for column S, row -> value (numbers) or $null (leave blank)
$values = @{
    3  = 2022
    4  = 10444.200000000001
    5  = 21.7
    6  = 7361.6
    7  = 143.1
    8  = 844.2
    9  = $null
    10 = "2 756,0"
    11 = "1 013,8"
    12 = "1 451,1"
    13 = 273.39999999999998
    14 = "-"
    15 = 17.7
}

for ($row = 3; $row -le 15; $row++) {
    $srcCell = $ws.Range("R$row")
    $dstCell = $ws.Range("S$row")

    # Copy the formatting (number format, font, borders, alignment, ...)
    # from the corresponding "2021" cell so the new column matches the
    # table's existing look.
    $srcCell.Copy()
    $dstCell.PasteSpecial($xlPasteFormats)

    $val = $values[$row]
    if ($null -ne $val) {
        $dstCell.Value = $val
    }
}

# Move the selection to T3, as in the source workbook after the edit.
$ws.Range("T3").Select()
